$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-unused trailing rows (91-100), working bottom-up so row
# indices of earlier rows are not disturbed mid-loop.
for ($r = 100; $r -ge 91; $r--) {
    $ws.Rows($r).Delete() | Out-Null
}

# Write the final Cluster name / Active cases table (header + 89 data rows).
$ws.Cells.Item(1, 1).Value = 'Cluster name'
$ws.Cells.Item(1, 2).Value = 'Active cases'
$ws.Cells.Item(2, 1).Value = '126 Racecourse Road Public Housing Tower Flemington'
$ws.Cells.Item(2, 2).Value = 7
$ws.Cells.Item(3, 1).Value = '3535 Opal Meadow Heights Aged Care Community Meadow Heights'
$ws.Cells.Item(3, 2).Value = 26
$ws.Cells.Item(4, 1).Value = 'Al Haj Halal Meats Glenroy'
$ws.Cells.Item(4, 2).Value = 11
$ws.Cells.Item(5, 1).Value = 'Al-Taqwa College Truganina'
$ws.Cells.Item(5, 2).Value = 8
$ws.Cells.Item(6, 1).Value = 'Apartment Complex 120 Racecourse Road North Melbourne'
$ws.Cells.Item(6, 2).Value = 5
$ws.Cells.Item(7, 1).Value = 'Australia Post Distribution Centre Sunshine West'
$ws.Cells.Item(7, 2).Value = 5
$ws.Cells.Item(8, 1).Value = 'Australian Lamb Colac East'
$ws.Cells.Item(8, 2).Value = 5
$ws.Cells.Item(9, 1).Value = 'Baxter Foods Australia Campbellfield'
$ws.Cells.Item(9, 2).Value = 5
$ws.Cells.Item(10, 1).Value = 'CFMEU Melbourne Office'
$ws.Cells.Item(10, 2).Value = 5
$ws.Cells.Item(11, 1).Value = 'CS Square Caroline Springs'
$ws.Cells.Item(11, 2).Value = 11
$ws.Cells.Item(12, 1).Value = 'Cafe Roco Dandenong'
$ws.Cells.Item(12, 2).Value = 6
$ws.Cells.Item(13, 1).Value = 'Campbellfield Ford Complex Vaccination Clinic Campbellfield'
$ws.Cells.Item(13, 2).Value = 8
$ws.Cells.Item(14, 1).Value = 'Cardinia Lakes Early Learning Centre Pakenham'
$ws.Cells.Item(14, 2).Value = 5
$ws.Cells.Item(15, 1).Value = 'Caroline Springs Police Station'
$ws.Cells.Item(15, 2).Value = 8
$ws.Cells.Item(16, 1).Value = 'Chemist Warehouse Campbellfield DC'
$ws.Cells.Item(16, 2).Value = 5
$ws.Cells.Item(17, 1).Value = 'Chemist Warehouse Fillo Drive Somerton'
$ws.Cells.Item(17, 2).Value = 12
$ws.Cells.Item(18, 1).Value = 'City of Wyndham Community'
$ws.Cells.Item(18, 2).Value = 6
$ws.Cells.Item(19, 1).Value = 'Coles Campbellfield Plaza Campbellfield'
$ws.Cells.Item(19, 2).Value = 7
$ws.Cells.Item(20, 1).Value = 'Coles Coburg North Village'
$ws.Cells.Item(20, 2).Value = 17
$ws.Cells.Item(21, 1).Value = 'Coles Pakenham Place Shopping Centre'
$ws.Cells.Item(21, 2).Value = 6
$ws.Cells.Item(22, 1).Value = 'Community Kids Bayswater Early Education Centre Bayswater North'
$ws.Cells.Item(22, 2).Value = 16
$ws.Cells.Item(23, 1).Value = 'Construction Site 1 Warde Street Footscray'
$ws.Cells.Item(23, 2).Value = 5
$ws.Cells.Item(24, 1).Value = 'Construction Site Olea Apartment Caulfield North'
$ws.Cells.Item(24, 2).Value = 14
$ws.Cells.Item(25, 1).Value = 'Costco Wholesale Epping'
$ws.Cells.Item(25, 2).Value = 28
$ws.Cells.Item(26, 1).Value = 'Crusader Caravans Epping'
$ws.Cells.Item(26, 2).Value = 22
$ws.Cells.Item(27, 1).Value = 'Dandenong Police Station Dandenong'
$ws.Cells.Item(27, 2).Value = 6
$ws.Cells.Item(28, 1).Value = 'DayHab Rehabilitation Treatment Centre Ringwood East'
$ws.Cells.Item(28, 2).Value = 6
$ws.Cells.Item(29, 1).Value = 'Direct Freight Express Campbellfield'
$ws.Cells.Item(29, 2).Value = 7
$ws.Cells.Item(30, 1).Value = 'Disability Residence Life without Barriers Ashwood'
$ws.Cells.Item(30, 2).Value = 5
$ws.Cells.Item(31, 1).Value = 'Don Watson Coldstore Derrimutg'
$ws.Cells.Item(31, 2).Value = 5
$ws.Cells.Item(32, 1).Value = 'Epworth Healthcare Epworth Richmond Emergency Department'
$ws.Cells.Item(32, 2).Value = 6
$ws.Cells.Item(33, 1).Value = 'Ermha365 Ltd Doveton'
$ws.Cells.Item(33, 2).Value = 9
$ws.Cells.Item(34, 1).Value = 'Fine Food Holdings Pty Ltd Dandenong South'
$ws.Cells.Item(34, 2).Value = 10
$ws.Cells.Item(35, 1).Value = 'Fonterra Manufacturing Workplace Campbellfield'
$ws.Cells.Item(35, 2).Value = 9
$ws.Cells.Item(36, 1).Value = 'General Foods Campbellfield'
$ws.Cells.Item(36, 2).Value = 9
$ws.Cells.Item(37, 1).Value = 'General Foods Campbellfield'
$ws.Cells.Item(37, 2).Value = 11
$ws.Cells.Item(38, 1).Value = 'Gladstone Parade Early Learning & Kinder Glenroy'
$ws.Cells.Item(38, 2).Value = 7
$ws.Cells.Item(39, 1).Value = 'Goodstart Early Learning Altona'
$ws.Cells.Item(39, 2).Value = 9
$ws.Cells.Item(40, 1).Value = 'Green Leaves Early Learning Cairnlea'
$ws.Cells.Item(40, 2).Value = 5
$ws.Cells.Item(41, 1).Value = 'Green Leaves Early Learning Centre Highlands Craigieburn'
$ws.Cells.Item(41, 2).Value = 16
$ws.Cells.Item(42, 1).Value = 'Greenvale Primary School'
$ws.Cells.Item(42, 2).Value = 5
$ws.Cells.Item(43, 1).Value = 'Hamilton Marino 236 Jasper Road McKinnon'
$ws.Cells.Item(43, 2).Value = 9
$ws.Cells.Item(44, 1).Value = 'ISS Factory Level 1 Terminal 2 Melbourne Airport Tullamarine'
$ws.Cells.Item(44, 2).Value = 10
$ws.Cells.Item(45, 1).Value = 'Ibis Kingsgate Hotel Melbourne'
$ws.Cells.Item(45, 2).Value = 6
$ws.Cells.Item(46, 1).Value = 'Ilim College Kiewa Campus Boys Secondary Dallas'
$ws.Cells.Item(46, 2).Value = 5
$ws.Cells.Item(47, 1).Value = 'Industrial Galvanizers Valmont Coatings Campbellfield'
$ws.Cells.Item(47, 2).Value = 14
$ws.Cells.Item(48, 1).Value = 'Inghams Enterprises Thomastown'
$ws.Cells.Item(48, 2).Value = 6
$ws.Cells.Item(49, 1).Value = 'Kippers Seafood Werribee'
$ws.Cells.Item(49, 2).Value = 6
$ws.Cells.Item(50, 1).Value = 'Kool Kidz Childcare Narre Warren'
$ws.Cells.Item(50, 2).Value = 13
$ws.Cells.Item(51, 1).Value = 'Level Crossing Removal Project Lilydale Construction Site John Street'
$ws.Cells.Item(51, 2).Value = 6
$ws.Cells.Item(52, 1).Value = 'Lineage Logistics Laverton North'
$ws.Cells.Item(52, 2).Value = 6
$ws.Cells.Item(53, 1).Value = 'Linfox Somerton National Distribution Centre Somerton'
$ws.Cells.Item(53, 2).Value = 10
$ws.Cells.Item(54, 1).Value = 'Mecca D.C Warehouse Melbourne Airport'
$ws.Cells.Item(54, 2).Value = 8
$ws.Cells.Item(55, 1).Value = 'Melbourne Assessment Prison West Melbourne'
$ws.Cells.Item(55, 2).Value = 7
$ws.Cells.Item(56, 1).Value = 'Melbourne Metropolitan Remand Centre Ravenhall'
$ws.Cells.Item(56, 2).Value = 8
$ws.Cells.Item(57, 1).Value = 'Melbourne West Police Station Docklands'
$ws.Cells.Item(57, 2).Value = 5
$ws.Cells.Item(58, 1).Value = 'Mill Park Police Station Mill Park'
$ws.Cells.Item(58, 2).Value = 8
$ws.Cells.Item(59, 1).Value = 'MyCentre Childcare Broadmeadows'
$ws.Cells.Item(59, 2).Value = 21
$ws.Cells.Item(60, 1).Value = 'National Gallery of Victoria Melbourne'
$ws.Cells.Item(60, 2).Value = 8
$ws.Cells.Item(61, 1).Value = 'Nido Early School Ascot Vale'
$ws.Cells.Item(61, 2).Value = 12
$ws.Cells.Item(62, 1).Value = 'Nido Early School Glenroy'
$ws.Cells.Item(62, 2).Value = 22
$ws.Cells.Item(63, 1).Value = 'Northern Health Northern Hospital Epping Emergency Department Tier 1B'
$ws.Cells.Item(63, 2).Value = 55
$ws.Cells.Item(64, 1).Value = 'Northern Health The Northern Hospital Epping'
$ws.Cells.Item(64, 2).Value = 16
$ws.Cells.Item(65, 1).Value = 'OnQ Plumbing and Excavations Craigieburn'
$ws.Cells.Item(65, 2).Value = 12
$ws.Cells.Item(66, 1).Value = 'Oporto Coolaroo'
$ws.Cells.Item(66, 2).Value = 8
$ws.Cells.Item(67, 1).Value = 'Oscar Romero Catholic Primary School Craigieburn'
$ws.Cells.Item(67, 2).Value = 5
$ws.Cells.Item(68, 1).Value = 'Our Lady Help of Christian''s Primary School Brunswick East'
$ws.Cells.Item(68, 2).Value = 9
$ws.Cells.Item(69, 1).Value = 'Pacific Meat Thomastown'
$ws.Cells.Item(69, 2).Value = 5
$ws.Cells.Item(70, 1).Value = 'Private Residence Daycare Allumba Way Wollert'
$ws.Cells.Item(70, 2).Value = 8
$ws.Cells.Item(71, 1).Value = 'Ramsay Health Care Warrigal Private Hospital'
$ws.Cells.Item(71, 2).Value = 5
$ws.Cells.Item(72, 1).Value = 'Ravenhall Correctional Centre Ravenhall'
$ws.Cells.Item(72, 2).Value = 8
$ws.Cells.Item(73, 1).Value = 'Richmond Quarter 261-271 Bridge Road Construction Site Richmond'
$ws.Cells.Item(73, 2).Value = 12
$ws.Cells.Item(74, 1).Value = 'Sacca''s Fruit World Broadmeadows Central Shopping Centre'
$ws.Cells.Item(74, 2).Value = 6
$ws.Cells.Item(75, 1).Value = 'Sharpline Stainless Steel Coburg North'
$ws.Cells.Item(75, 2).Value = 6
$ws.Cells.Item(76, 1).Value = 'St Margaret''s Primary School OSHC Maribyrnong'
$ws.Cells.Item(76, 2).Value = 12
$ws.Cells.Item(77, 1).Value = 'St Vincents Hospital Emergency Department Melbourne'
$ws.Cells.Item(77, 2).Value = 9
$ws.Cells.Item(78, 1).Value = 'Tek Foods Somerton'
$ws.Cells.Item(78, 2).Value = 10
$ws.Cells.Item(79, 1).Value = 'The Huntly-Goornong Rail Works'
$ws.Cells.Item(79, 2).Value = 6
$ws.Cells.Item(80, 1).Value = 'The Royal Children''s Hospital Melbourne Emergency Department Parkville Tier 1B'
$ws.Cells.Item(80, 2).Value = 18
$ws.Cells.Item(81, 1).Value = 'The Royal Melbourne Hospital Ward 6SE Parkville'
$ws.Cells.Item(81, 2).Value = 21
$ws.Cells.Item(82, 1).Value = 'The Royal Talbot Rehabilitation Centre Kew'
$ws.Cells.Item(82, 2).Value = 10
$ws.Cells.Item(83, 1).Value = 'ThorwestenCabinets Pakenham'
$ws.Cells.Item(83, 2).Value = 14
$ws.Cells.Item(84, 1).Value = 'Truganina Early Learning Centre Truganina'
$ws.Cells.Item(84, 2).Value = 6
$ws.Cells.Item(85, 1).Value = 'Unilodge College Square Student Accommodation 570 Lygon Street Carlton'
$ws.Cells.Item(85, 2).Value = 6
$ws.Cells.Item(86, 1).Value = 'Wallaby Childcare Wollert'
$ws.Cells.Item(86, 2).Value = 18
$ws.Cells.Item(87, 1).Value = 'Werribee Mercy Hospital Emergency Department'
$ws.Cells.Item(87, 2).Value = 13
$ws.Cells.Item(88, 1).Value = 'Western Health Footscray Hospital Emergency Department'
$ws.Cells.Item(88, 2).Value = 8
$ws.Cells.Item(89, 1).Value = 'Western Health Sunshine Hospital Emergency Department'
$ws.Cells.Item(89, 2).Value = 9
$ws.Cells.Item(90, 1).Value = 'Yarra Childcare Centre Truganina'
$ws.Cells.Item(90, 2).Value = 5
